# Reverses the row order of the "Periodo Mora" / "Valor Mora" table
# (E16:F48) on Hoja1: the table previously listed periods 2111..2407
# in ascending order with a single differing "Valor Mora" on the last
# period (2407). The updated database now lists the same periods in
# descending order (2407..2111), so the differing value now lands on
# the first row instead of the last.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow = 48
$rowCount = $lastRow - $firstRow + 1

# Capture the current (pre-edit) values for the period codes (E) and
# the mora amounts (F) before overwriting anything.
$periods = @()
$amounts = @()
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $periods += $ws.Cells.Item($r, 5).Value2
    $amounts += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse row order, leaving every other column
# (and all cell formatting) untouched.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $srcIndex = $rowCount - 1 - $i
    $ws.Cells.Item($r, 5).Value = $periods[$srcIndex]
    $ws.Cells.Item($r, 6).Value = $amounts[$srcIndex]
}
